$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 30, shifting existing rows 30-33 down to 31-34.
$ws.Rows.Item(30).Insert()

# Populate the new row 30 with a copy of the data that used to live in row 30
# (same market/product info), but with an updated date (Fecha).
$ws.Cells.Item(30, 1).Value = 6
$ws.Cells.Item(30, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(30, 3).Value = "Metropolitana"
$ws.Cells.Item(30, 4).Value = 44641
$ws.Cells.Item(30, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(30, 5).Value = 13
$ws.Cells.Item(30, 6).Value = "Fruta"
$ws.Cells.Item(30, 7).Value = 100102
$ws.Cells.Item(30, 8).Value = "Cítricos"
$ws.Cells.Item(30, 9).Value = 100102006
$ws.Cells.Item(30, 10).Value = "Pomelo"
$ws.Cells.Item(30, 11).Value = "Start Ruby"
$ws.Cells.Item(30, 12).Value = "Primera"
$ws.Cells.Item(30, 13).Value = 16
$ws.Cells.Item(30, 14).Value = 180000
$ws.Cells.Item(30, 15).Value = 180000
$ws.Cells.Item(30, 16).Value = 180000
$ws.Cells.Item(30, 17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(30, 18).Value = "Región Metropolitana"
$ws.Cells.Item(30, 19).Value = 514
$ws.Cells.Item(30, 20).Value = 350
